$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 108 (shifts old rows
# 108-119 down to 110-121), mirroring the weekly-update pattern where the
# newest data point is added at the top of this market/product block.
$ws.Range("A108:A109").EntireRow.Insert()

# New row 108: Primera quality, new weekly observation (2022-12-23)
$ws.Cells.Item(108, 1).Value = 11
$ws.Cells.Item(108, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(108, 3).Value = "Bíobío"
$ws.Cells.Item(108, 4).Value = 44918
$ws.Cells.Item(108, 5).Value = 8
$ws.Cells.Item(108, 6).Value = "Fruta"
$ws.Cells.Item(108, 7).Value = 100101
$ws.Cells.Item(108, 8).Value = "Berries"
$ws.Cells.Item(108, 9).Value = 100101001
$ws.Cells.Item(108, 10).Value = "Arándano (blue)"
$ws.Cells.Item(108, 11).Value = "Sin especificar"
$ws.Cells.Item(108, 12).Value = "Primera"
$ws.Cells.Item(108, 13).Value = 200
$ws.Cells.Item(108, 14).Value = 3000
$ws.Cells.Item(108, 15).Value = 3500
$ws.Cells.Item(108, 16).Value = 3250
$ws.Cells.Item(108, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(108, 18).Value = "Región de Ñuble"
$ws.Cells.Item(108, 19).Value = 1625
$ws.Cells.Item(108, 20).Value = 2

# New row 109: Segunda quality, same new weekly observation (2022-12-23)
$ws.Cells.Item(109, 1).Value = 11
$ws.Cells.Item(109, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(109, 3).Value = "Bíobío"
$ws.Cells.Item(109, 4).Value = 44918
$ws.Cells.Item(109, 5).Value = 8
$ws.Cells.Item(109, 6).Value = "Fruta"
$ws.Cells.Item(109, 7).Value = 100101
$ws.Cells.Item(109, 8).Value = "Berries"
$ws.Cells.Item(109, 9).Value = 100101001
$ws.Cells.Item(109, 10).Value = "Arándano (blue)"
$ws.Cells.Item(109, 11).Value = "Sin especificar"
$ws.Cells.Item(109, 12).Value = "Segunda"
$ws.Cells.Item(109, 13).Value = 100
$ws.Cells.Item(109, 14).Value = 2500
$ws.Cells.Item(109, 15).Value = 2500
$ws.Cells.Item(109, 16).Value = 2500
$ws.Cells.Item(109, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(109, 18).Value = "Región de Ñuble"
$ws.Cells.Item(109, 19).Value = 1250
$ws.Cells.Item(109, 20).Value = 2
